$wb = $excel.ActiveWorkbook

# --- 1) Global rename: "MF3H" -> "T3H" wherever it survives in the final data ---
$wsEq = $wb.Worksheets.Item("equilibrium_concentrations")
$wsEq.Range("B1").Value = "T3H"

$wsStoich = $wb.Worksheets.Item("input_stoich_coefficients")
$wsStoich.Range("B1").Value = "T3H"

$wsEnth = $wb.Worksheets.Item("input_enthalpies")
$wsEnth.Range("A3").Value = "T3H"

# --- 2) input_enthalpies: clear row 2 (A2:B2) ---
$wsEnth.Range("A2:B2").ClearContents()

# --- 3) constants_evaluated: new headers, single remaining data row ---
$wsConst = $wb.Worksheets.Item("constants_evaluated")
$wsConst.Range("A1").Value = "Component"
$wsConst.Range("B1").Value = "Constant"
$wsConst.Range("C1").Value = "St.Deviation"
$wsConst.Range("D1").Value = "Validity"
$wsConst.Range("A2").Value = "Comp"
$wsConst.Range("B2:D2").NumberFormat = "@"
$wsConst.Range("B2").Value = "5.4061110496521"
$wsConst.Range("C2").Value = "0.0205130354010432"
$wsConst.Range("D2").Value = "-Inf"
$wsConst.Range("A3:D4").ClearContents()

# --- 4) enthalpies_calculated: single remaining data row ---
$wsEnthCalc = $wb.Worksheets.Item("enthalpies_calculated")
$wsEnthCalc.Range("A2").Value = "Comp"
$wsEnthCalc.Range("B2").Value = 50.7744283493441
$wsEnthCalc.Range("C2").Value = 0.534436612020058
$wsEnthCalc.Range("A3:C4").ClearContents()

# --- 5) input_stoich_coefficients: single remaining data row ---
$wsStoich.Range("A2:B2").NumberFormat = "@"
$wsStoich.Range("A2").Value = "1"
$wsStoich.Range("B2").Value = "1"
$wsStoich.Range("C2").Value = "Comp"
$wsStoich.Range("A3:C4").ClearContents()
